$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 42 (Peach Saison): add a Miscellaneous Notes entry
$ws.Range("M42").Value = "Aged in chardonnay barrels"

# Row 44: replace the old "Monks' IPA / NEIPA" entry with the new
# Spencer seasonal beer "The Monkster Mash" (Pumpkin Ale)
$ws.Range("C44").Value = "The Monkster Mash"
$ws.Range("D44").Value = "Pumpkin Ale"
$ws.Range("E44").Value = 5.2
$ws.Range("H44").Value = "Citra"
$ws.Range("I44").Value = "pumpkin, spice"
$ws.Range("M44").Value = "Seasonal, canned"
